$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header row: "<Label>_old" -> "<Label>_FV2404" (A1:J1)
# and "<Label>_new" -> "<Label>_FV2410" (L1:U1). K1 ("diff") is unchanged.
$labels = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $labels.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($labels[$i])_FV2404"
    $ws.Cells.Item(1, $i + 12).Value = "$($labels[$i])_FV2410"
}

# Turn the A1:U79 range into an Excel Table (ListObject) with an AutoFilter,
# picking up the just-renamed header row as the column names.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U79"), $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (split below row 1).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
